$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C) for rows 2-7 from 45170 to 45174
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
